# repull data, push all data, mean calculation
# Update column F (dSF) values for the affected rows to reflect
# repulled data used in the mean calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12 = 0
    14 = -2
    20 = -4
    24 = 0
    26 = -2
    28 = 0
    31 = -2
    35 = 0
    37 = -3
    40 = 2
    42 = -4
    43 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
